$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column R (18th column) to make room for "Periodicity"
$ws.Columns.Item(18).Insert()

# Header for the new column
$ws.Cells.Item(1, 18).Value = "Periodicity"

# Fill "Daily" for rows 2 through 21 in the new column R
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 18).Value = "Daily"
}
